# Update the registered e-mail addresses on the "Teste" sheet.
# peterpan17@movie.com -> peterpan20@movie.com (D3)
# peterpan18@movie.com -> peterpan21@movie.com (D4)
#
# Both cells already carry Excel's "quote prefix" (stored-as-text) cell
# format, the same one used on B3/B4 (the hyperlinked site column). A
# leading apostrophe keeps that formatting intact when the new value is
# typed in, instead of Excel silently re-styling the cell as it would if
# the apostrophe were omitted.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Teste")

$ws.Range("D3").Value = "'peterpan20@movie.com"
$ws.Range("D4").Value = "'peterpan21@movie.com"
